$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44468
$ws.Cells.Item(2, 14).Value = 29000
$ws.Cells.Item(2, 15).Value = 30000
$ws.Cells.Item(2, 16).Value = 29500
$ws.Cells.Item(2, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(2, 19).Value = 2950
$ws.Cells.Item(2, 20).Value = 10

# Row 3
$ws.Cells.Item(3, 4).Value = 44475
$ws.Cells.Item(3, 12).Value = 'Especial'
$ws.Cells.Item(3, 13).Value = 200
$ws.Cells.Item(3, 14).Value = 32000
$ws.Cells.Item(3, 15).Value = 33000
$ws.Cells.Item(3, 16).Value = 32500
$ws.Cells.Item(3, 19).Value = 2708

# Row 4
$ws.Cells.Item(4, 4).Value = 44167
$ws.Cells.Item(4, 13).Value = 200
$ws.Cells.Item(4, 14).Value = 18000
$ws.Cells.Item(4, 15).Value = 19000
$ws.Cells.Item(4, 16).Value = 18500
$ws.Cells.Item(4, 17).Value = '$/caja 13 kilos'
$ws.Cells.Item(4, 19).Value = 1423
$ws.Cells.Item(4, 20).Value = 13

# Row 5
$ws.Cells.Item(5, 4).Value = 44496
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 13).Value = 200
$ws.Cells.Item(5, 14).Value = 23000
$ws.Cells.Item(5, 15).Value = 24000
$ws.Cells.Item(5, 16).Value = 23500
$ws.Cells.Item(5, 19).Value = 1958

# Row 6
$ws.Cells.Item(6, 4).Value = 44839
$ws.Cells.Item(6, 12).Value = 'Segunda'
$ws.Cells.Item(6, 13).Value = 160
$ws.Cells.Item(6, 14).Value = 26000
$ws.Cells.Item(6, 15).Value = 27000
$ws.Cells.Item(6, 16).Value = 26500
$ws.Cells.Item(6, 17).Value = '$/caja 12 kilos'
$ws.Cells.Item(6, 19).Value = 2208
$ws.Cells.Item(6, 20).Value = 12

# Row 7
$ws.Cells.Item(7, 4).Value = 44881
$ws.Cells.Item(7, 14).Value = 22000
$ws.Cells.Item(7, 15).Value = 23000
$ws.Cells.Item(7, 16).Value = 22500
$ws.Cells.Item(7, 17).Value = '$/caja 12 kilos'
$ws.Cells.Item(7, 19).Value = 1875

# Row 8
$ws.Cells.Item(8, 4).Value = 44874
$ws.Cells.Item(8, 12).Value = 'Segunda'
$ws.Cells.Item(8, 13).Value = 250
$ws.Cells.Item(8, 14).Value = 22000
$ws.Cells.Item(8, 15).Value = 23000
$ws.Cells.Item(8, 16).Value = 22500
$ws.Cells.Item(8, 19).Value = 1875

# Row 9
$ws.Cells.Item(9, 4).Value = 45126

# Row 10
$ws.Cells.Item(10, 4).Value = 45126

# Row 11
$ws.Cells.Item(11, 4).Value = 45125
$ws.Cells.Item(11, 13).Value = 160
$ws.Cells.Item(11, 14).Value = 14000
$ws.Cells.Item(11, 15).Value = 15000
$ws.Cells.Item(11, 16).Value = 14375
$ws.Cells.Item(11, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(11, 19).Value = 1438
$ws.Cells.Item(11, 20).Value = 10

# Row 12
$ws.Cells.Item(12, 4).Value = 45125
$ws.Cells.Item(12, 13).Value = 180
$ws.Cells.Item(12, 14).Value = 13000
$ws.Cells.Item(12, 15).Value = 13000
$ws.Cells.Item(12, 16).Value = 13000
$ws.Cells.Item(12, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(12, 19).Value = 1300

# Row 13
$ws.Cells.Item(13, 4).Value = 44783
$ws.Cells.Item(13, 12).Value = 'Tercera'
$ws.Cells.Item(13, 13).Value = 100
$ws.Cells.Item(13, 14).Value = 27000
$ws.Cells.Item(13, 15).Value = 28000
$ws.Cells.Item(13, 16).Value = 27500
$ws.Cells.Item(13, 19).Value = 2292

# Row 14
$ws.Cells.Item(14, 4).Value = 44545
$ws.Cells.Item(14, 12).Value = 'Primera'
$ws.Cells.Item(14, 13).Value = 200
$ws.Cells.Item(14, 14).Value = 23000
$ws.Cells.Item(14, 15).Value = 24000
$ws.Cells.Item(14, 16).Value = 23500
$ws.Cells.Item(14, 17).Value = '$/bandeja 12 kilos'
$ws.Cells.Item(14, 19).Value = 1958

# Row 15
$ws.Cells.Item(15, 4).Value = 44441
$ws.Cells.Item(15, 12).Value = 'Primera'
$ws.Cells.Item(15, 13).Value = 100
$ws.Cells.Item(15, 14).Value = 29000
$ws.Cells.Item(15, 15).Value = 30000
$ws.Cells.Item(15, 16).Value = 29500
$ws.Cells.Item(15, 19).Value = 2458

# Row 16
$ws.Cells.Item(16, 4).Value = 44776
$ws.Cells.Item(16, 12).Value = 'Segunda'
$ws.Cells.Item(16, 13).Value = 160
$ws.Cells.Item(16, 14).Value = 29000
$ws.Cells.Item(16, 15).Value = 30000
$ws.Cells.Item(16, 16).Value = 29500
$ws.Cells.Item(16, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(16, 18).Value = 'Región de Coquimbo'
$ws.Cells.Item(16, 19).Value = 2950
$ws.Cells.Item(16, 20).Value = 10

# Row 17
$ws.Cells.Item(17, 4).Value = 44160
$ws.Cells.Item(17, 12).Value = 'Segunda'
$ws.Cells.Item(17, 13).Value = 200
$ws.Cells.Item(17, 14).Value = 19000
$ws.Cells.Item(17, 15).Value = 20000
$ws.Cells.Item(17, 16).Value = 19500
$ws.Cells.Item(17, 17).Value = '$/caja 13 kilos'
$ws.Cells.Item(17, 19).Value = 1500
$ws.Cells.Item(17, 20).Value = 13

# Row 18
$ws.Cells.Item(18, 4).Value = 44811
$ws.Cells.Item(18, 12).Value = 'Primera'
$ws.Cells.Item(18, 14).Value = 29000
$ws.Cells.Item(18, 15).Value = 30000
$ws.Cells.Item(18, 16).Value = 29500
$ws.Cells.Item(18, 19).Value = 2458

# Row 19
$ws.Cells.Item(19, 4).Value = 44482
$ws.Cells.Item(19, 14).Value = 25000
$ws.Cells.Item(19, 15).Value = 26000
$ws.Cells.Item(19, 16).Value = 25500
$ws.Cells.Item(19, 17).Value = '$/caja 12 kilos'
$ws.Cells.Item(19, 19).Value = 2125
$ws.Cells.Item(19, 20).Value = 12

# Row 20
$ws.Cells.Item(20, 4).Value = 44489
$ws.Cells.Item(20, 12).Value = 'Primera'
$ws.Cells.Item(20, 13).Value = 200
$ws.Cells.Item(20, 14).Value = 24000
$ws.Cells.Item(20, 15).Value = 25000
$ws.Cells.Item(20, 16).Value = 24500
$ws.Cells.Item(20, 17).Value = '$/caja 12 kilos'
$ws.Cells.Item(20, 19).Value = 2042
$ws.Cells.Item(20, 20).Value = 12

# Row 21
$ws.Cells.Item(21, 4).Value = 44860
$ws.Cells.Item(21, 12).Value = 'Primera'
$ws.Cells.Item(21, 14).Value = 23000
$ws.Cells.Item(21, 15).Value = 24000
$ws.Cells.Item(21, 16).Value = 23500
$ws.Cells.Item(21, 17).Value = '$/caja 12 kilos'
$ws.Cells.Item(21, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(21, 19).Value = 1958
$ws.Cells.Item(21, 20).Value = 12

# Row 22
$ws.Cells.Item(22, 4).Value = 44524
$ws.Cells.Item(22, 13).Value = 200
$ws.Cells.Item(22, 14).Value = 23000
$ws.Cells.Item(22, 15).Value = 24000
$ws.Cells.Item(22, 16).Value = 23500
$ws.Cells.Item(22, 19).Value = 1958

# Row 23
$ws.Cells.Item(23, 4).Value = 44846
$ws.Cells.Item(23, 13).Value = 160
$ws.Cells.Item(23, 14).Value = 24000
$ws.Cells.Item(23, 15).Value = 25000
$ws.Cells.Item(23, 16).Value = 24500
$ws.Cells.Item(23, 17).Value = '$/caja 12 kilos'
$ws.Cells.Item(23, 19).Value = 2042
$ws.Cells.Item(23, 20).Value = 12

# Row 24
$ws.Cells.Item(24, 4).Value = 44846
$ws.Cells.Item(24, 12).Value = 'Segunda'
$ws.Cells.Item(24, 14).Value = 22000
$ws.Cells.Item(24, 15).Value = 23000
$ws.Cells.Item(24, 16).Value = 22500
$ws.Cells.Item(24, 19).Value = 1875

# Row 25
$ws.Cells.Item(25, 4).Value = 45133
$ws.Cells.Item(25, 12).Value = 'Primera'
$ws.Cells.Item(25, 13).Value = 150
$ws.Cells.Item(25, 14).Value = 15000
$ws.Cells.Item(25, 15).Value = 16000
$ws.Cells.Item(25, 16).Value = 15667
$ws.Cells.Item(25, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(25, 19).Value = 1567
$ws.Cells.Item(25, 20).Value = 10

Write-Host "Applied 169 cell updates"